$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated cryptos list values (price and volume-change columns).
# Numeric-looking price strings are written via a temporary Text
# number format so Excel keeps them as literal strings (e.g. "1.002")
# instead of auto-converting them to numbers; the style is reset
# back to Normal immediately afterwards so no visible formatting
# change is left behind.
# Also includes a rank swap between HuobiToken (row 32) and Filecoin (row 33).
$ws.Range("D2").Value = "27.644.07"
$ws.Range("E2").Value = "  -2.35%  "
$ws.Range("D3").Value = "1.761.08"
$ws.Range("E3").Value = "  -3.18%  "
$cell = $ws.Range("D4")
$cell.NumberFormat = "@"
$cell.Value = "1.002"
$cell.Style = "Normal"
$ws.Range("E4").Value = "  -0.08%  "
$cell = $ws.Range("D5")
$cell.NumberFormat = "@"
$cell.Value = "325.09"
$cell.Style = "Normal"
$ws.Range("E5").Value = "  -0.95%  "
$ws.Range("E6").Value = "  -0.05%  "
$cell = $ws.Range("D7")
$cell.NumberFormat = "@"
$cell.Value = "0.4324"
$cell.Style = "Normal"
$ws.Range("E7").Value = "  -0.53%  "
$cell = $ws.Range("D8")
$cell.NumberFormat = "@"
$cell.Value = "0.3612"
$cell.Style = "Normal"
$ws.Range("E8").Value = "  -1.67%  "
$cell = $ws.Range("D9")
$cell.NumberFormat = "@"
$cell.Value = "0.07585"
$cell.Style = "Normal"
$ws.Range("E9").Value = "  -1.46%  "
$cell = $ws.Range("D10")
$cell.NumberFormat = "@"
$cell.Value = "42.16"
$cell.Style = "Normal"
$ws.Range("E10").Value = "  -6.28%  "
$cell = $ws.Range("D11")
$cell.NumberFormat = "@"
$cell.Value = "1.113"
$cell.Style = "Normal"
$ws.Range("E11").Value = "  -2.85%  "
$cell = $ws.Range("D12")
$cell.NumberFormat = "@"
$cell.Value = "1.002"
$cell.Style = "Normal"
$ws.Range("E12").Value = "  -0.02%  "
$cell = $ws.Range("D13")
$cell.NumberFormat = "@"
$cell.Value = "20.84"
$cell.Style = "Normal"
$ws.Range("E13").Value = "  -5.93%  "
$cell = $ws.Range("D14")
$cell.NumberFormat = "@"
$cell.Value = "6.080"
$cell.Style = "Normal"
$ws.Range("E14").Value = "  -3.62%  "
$cell = $ws.Range("D15")
$cell.NumberFormat = "@"
$cell.Value = "7.228"
$cell.Style = "Normal"
$ws.Range("E15").Value = "  -4.17%  "
$ws.Range("D16").Value = "1.756.87"
$ws.Range("E16").Value = "  -4.39%  "
$cell = $ws.Range("D17")
$cell.NumberFormat = "@"
$cell.Value = "92.48"
$cell.Style = "Normal"
$ws.Range("E17").Value = "  -1.34%  "
$cell = $ws.Range("D18")
$cell.NumberFormat = "@"
$cell.Value = "0.00001070"
$cell.Style = "Normal"
$ws.Range("E18").Value = "  -1.27%  "
$cell = $ws.Range("D19")
$cell.NumberFormat = "@"
$cell.Value = "0.06424"
$cell.Style = "Normal"
$ws.Range("E19").Value = "  -2.29%  "
$ws.Range("E20").Value = "  +0.00%  "
$cell = $ws.Range("D21")
$cell.NumberFormat = "@"
$cell.Value = "17.12"
$cell.Style = "Normal"
$ws.Range("E21").Value = "  -2.43%  "
$cell = $ws.Range("D22")
$cell.NumberFormat = "@"
$cell.Value = "5.864"
$cell.Style = "Normal"
$ws.Range("D23").Value = "27.683.37"
$ws.Range("E23").Value = "  -2.31%  "
$cell = $ws.Range("D24")
$cell.NumberFormat = "@"
$cell.Value = "11.30"
$cell.Style = "Normal"
$ws.Range("E24").Value = "  -3.07%  "
$cell = $ws.Range("D25")
$cell.NumberFormat = "@"
$cell.Value = "2.106"
$cell.Style = "Normal"
$ws.Range("E25").Value = "  +2.02%  "
$cell = $ws.Range("D26")
$cell.NumberFormat = "@"
$cell.Value = "162.28"
$cell.Style = "Normal"
$ws.Range("E26").Value = "  -0.30%  "
$cell = $ws.Range("D27")
$cell.NumberFormat = "@"
$cell.Value = "20.53"
$cell.Style = "Normal"
$ws.Range("E27").Value = "  -1.01%  "
$ws.Range("D28").Value = "1.959.10"
$ws.Range("E28").Value = "  -4.01%  "
$ws.Range("E29").Value = "  -6.45%  "
$cell = $ws.Range("D30")
$cell.NumberFormat = "@"
$cell.Value = "126.13"
$cell.Style = "Normal"
$ws.Range("E30").Value = "  -2.31%  "
$cell = $ws.Range("D31")
$cell.NumberFormat = "@"
$cell.Value = "1.105"
$cell.Style = "Normal"
$ws.Range("E31").Value = "  -10.46%  "
$ws.Range("B32").Value = "Filecoin"
$ws.Range("C32").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$cell = $ws.Range("D32")
$cell.NumberFormat = "@"
$cell.Value = "5.620"
$cell.Style = "Normal"
$ws.Range("E32").Value = "  -5.97%  "
$ws.Range("B33").Value = "HuobiToken"
$ws.Range("C33").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$cell = $ws.Range("D33")
$cell.NumberFormat = "@"
$cell.Value = "3.669"
$cell.Style = "Normal"
$ws.Range("E33").Value = "  +5.72%  "
$cell = $ws.Range("D34")
$cell.NumberFormat = "@"
$cell.Value = "0.08971"
$cell.Style = "Normal"
$ws.Range("E34").Value = "  -2.44%  "
$cell = $ws.Range("D35")
$cell.NumberFormat = "@"
$cell.Value = "12.26"
$cell.Style = "Normal"
$ws.Range("E35").Value = "  -5.80%  "
$cell = $ws.Range("D36")
$cell.NumberFormat = "@"
$cell.Value = "0.02308"
$cell.Style = "Normal"
$ws.Range("E36").Value = "  -2.12%  "
$cell = $ws.Range("D37")
$cell.NumberFormat = "@"
$cell.Value = "0.2118"
$cell.Style = "Normal"
$ws.Range("E37").Value = "  -2.85%  "
$cell = $ws.Range("D38")
$cell.NumberFormat = "@"
$cell.Value = "0.6421"
$cell.Style = "Normal"
$ws.Range("E38").Value = "  -2.50%  "
$cell = $ws.Range("D39")
$cell.NumberFormat = "@"
$cell.Value = "0.06024"
$cell.Style = "Normal"
$ws.Range("E39").Value = "  -2.98%  "
$cell = $ws.Range("D40")
$cell.NumberFormat = "@"
$cell.Value = "4.965"
$cell.Style = "Normal"
$ws.Range("E40").Value = "  -4.76%  "
$ws.Range("E41").Value = "  -0.71%  "
$ws.Range("E42").Value = "  +0.02%  "
$cell = $ws.Range("D43")
$cell.NumberFormat = "@"
$cell.Value = "1.398"
$cell.Style = "Normal"
$ws.Range("E43").Value = "  -3.03%  "
$cell = $ws.Range("D44")
$cell.NumberFormat = "@"
$cell.Value = "7.939"
$cell.Style = "Normal"
$ws.Range("E44").Value = "  -2.56%  "
$cell = $ws.Range("D45")
$cell.NumberFormat = "@"
$cell.Value = "13.42"
$cell.Style = "Normal"
$ws.Range("E45").Value = "  -3.76%  "
$cell = $ws.Range("D46")
$cell.NumberFormat = "@"
$cell.Value = "0.5942"
$cell.Style = "Normal"
$ws.Range("E46").Value = "  -2.96%  "
$cell = $ws.Range("D47")
$cell.NumberFormat = "@"
$cell.Value = "3.721"
$cell.Style = "Normal"
$ws.Range("E47").Value = "  -1.00%  "
$cell = $ws.Range("D48")
$cell.NumberFormat = "@"
$cell.Value = "1.991"
$cell.Style = "Normal"
$ws.Range("E48").Value = "  -1.92%  "
$cell = $ws.Range("D49")
$cell.NumberFormat = "@"
$cell.Value = "122.50"
$cell.Style = "Normal"
$ws.Range("E49").Value = "  -2.80%  "
$cell = $ws.Range("D50")
$cell.NumberFormat = "@"
$cell.Value = "1.171"
$cell.Style = "Normal"
$ws.Range("E50").Value = "  +0.85%  "
$cell = $ws.Range("D51")
$cell.NumberFormat = "@"
$cell.Value = "0.06882"
$cell.Style = "Normal"
$ws.Range("E51").Value = "  -1.83%  "
